# Apply updated "dSF" (column F) values on Sheet1, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -2
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = -1
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 4
$ws.Range("F28").Value = -1
$ws.Range("F34").Value = 0
$ws.Range("F40").Value = -1
$ws.Range("F42").Value = 2
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = -1
$ws.Range("F47").Value = -1
$ws.Range("F49").Value = 3
$ws.Range("F56").Value = 0
$ws.Range("F57").Value = -2
$ws.Range("F63").Value = -3
$ws.Range("F64").Value = 1
$ws.Range("F73").Value = -4
